$d = $word.ActiveDocument

# 1. Merge "Sigurd Skov Jensen" split runs into one run (just normalize text via Find/Replace)
$d.Content.Find.Execute("Sigurd Skov Jensen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sigurd Skov Jensen", 2)

# 2. Add a new row to the table with Jeppe Qvistgaard Hansen / 201703746
$table = $d.Tables.Item(1)
$newRow = $table.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Jeppe Qvistgaard Hansen"
$newRow.Cells.Item(2).Range.Text = "201703746"

# 3. Merge "URL for github: " split runs into one run
$d.Content.Find.Execute("URL for github: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "URL for github: ", 2)
